$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 3 and 4 content ---
# New Row 3 (previously row 4's data), with explicit black font style (s="2")
$row3 = @("Howe","Andrew",16,"10th","Male","USA/Philippines ",20170429,150154,
  "This passage is adapted from Francis J. Flynn and Gabrielle S. Adams, `"Money Can't Buy Love: Asymmetric Beliefs about Gift Price and Feelings of Appreciation.`" ©2008 by Elsevier Inc.",
  7,4,0,11)
for ($col = 1; $col -le 13; $col++) {
  $c = $ws.Cells.Item(3, $col)
  $c.Value = $row3[$col - 1]
  $c.Font.Color = 0
}

# New Row 4 (previously row 3's data), default style (no explicit font color)
$row4 = @("Kraft","Justin",15,"10th","Male","USA/Japanese",20170430,144237,
  "This passage is adapted from Saki, `“The Schartz-Metterklume Method.`” Originally published in 1911.",
  7,3,0,10)
for ($col = 1; $col -le 13; $col++) {
  $c = $ws.Cells.Item(4, $col)
  $c.Value = $row4[$col - 1]
}
$ws.Range("A4:M4").ClearFormats()

# --- Append new rows 13-19 with new test subject info / scores ---
# Row data, keyed by target row number
$newRows = @{
  13 = @("Ambrosino ","Jack",17,"11th","Male","USA",20170617,170926,"Test 8 pas 1",10,0,0,10)
  14 = @("Ambrosino ","Jack",17,"11th","Male","USA",20170617,192567,"Test 8 pas 2",10,1,0,10)
  15 = @("Black","Nicholas",17,"11th","Male","USA/Japanese",20170617,174536,"test 1 pas 1",9,1,0,10)
  16 = @("Black","Nicholas",17,"11th","Male","USA/Japanese",20170617,175545,"tes 1 pas 2",10,1,0,11)
  17 = @("Howe","Alexis",17,"11th","Female","USA/Philippines ",20170619,174323,"TEST 2 passage 1 ",10,0,0,10)
  18 = @("Howe","Alexis",17,"11th","Female ","USA/Philippines ",20170619,175643,"TEST 2 passage 2",11,0,0,11)
  19 = @("Sasanuma","Chris",16,"10th","Male","USA/Japanese",20170624,174536,"test 3 pas 1",8,2,0,10)
}

# Original entry order (matches how the shared-strings table ended up ordered):
# row 13, then 15, 16, 17, 18, then back to 14, then 19.
$entryOrder = @(13, 15, 16, 17, 18, 14, 19)

foreach ($r in $entryOrder) {
  $rowData = $newRows[$r]
  for ($col = 1; $col -le 13; $col++) {
    $ws.Cells.Item($r, $col).Value = $rowData[$col - 1]
  }
}

# --- Update the active cell selection to reflect the new extent ---
$ws.Range("H20").Select()

Write-Host "done"
